# Commit: "adding new progress as of date 04 nov 2025"
#
# For rows 3..29 on the "Training Dashboard" sheet:
#   - column H (PERIOD TO EXPIRE) decreases by 1 (one day closer to expiry)
#   - column I (LAST UPDATE) changes from "03-Nov-2025" to "04-Nov-2025"
#
# Note: assigning a date-looking literal string straight into Range.Value
# (or Value2) on a "General" formatted cell makes Excel auto-convert it into
# a real date serial number and silently reformats the cell's NumberFormat,
# which is not what the source workbook does (column I stores a plain text
# value). To keep the cell as literal text -- with its original style/
# NumberFormat untouched -- the new text is produced via a text formula and
# then solidified in place with a Copy / PasteSpecial (values only), which
# does not go through Excel's "looks like a date" auto-detection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$xlPasteValues = -4163

for ($row = 3; $row -le 29; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H - PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($row, 9)   # column I - LAST UPDATE

    $hCell.Value = $hCell.Value2 - 1

    $iCell.Formula = '="04-Nov-2025"'
    $iCell.Copy()
    $iCell.PasteSpecial($xlPasteValues)
}

$excel.CutCopyMode = 0
